$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3, shifting old row 3 (and below) down
$ws.Rows.Item(3).Insert()

# New row 3: copy the values from the row below (old row 3, now row 4) but use newer date & prices from row 2
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(3, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(3, 4).Value = 44533
$ws.Cells.Item(3, 4).NumberFormat = $ws.Cells.Item(4, 4).NumberFormat
$ws.Cells.Item(3, 5).Value = 15
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100108
$ws.Cells.Item(3, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(3, 9).Value = 100108007
$ws.Cells.Item(3, 10).Value = "Coco"
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 100
$ws.Cells.Item(3, 14).Value = 16000
$ws.Cells.Item(3, 15).Value = 17000
$ws.Cells.Item(3, 16).Value = 16500
$ws.Cells.Item(3, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(3, 18).Value = "Perú"
$ws.Cells.Item(3, 19).Value = 825
$ws.Cells.Item(3, 20).Value = 20
